$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape (GitHub Actions run)

$ws.Range("D2").Value = "'58.276.32"
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").Value = "'2.595.08"
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'522.30"
$ws.Range("E5").Value = '  +1.14%  '

$ws.Range("D6").Value = "'144.60"
$ws.Range("E6").Value = '  +1.78%  '

$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").Value = "'2.614.72"
$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").Value = "'6.66"
$ws.Range("E10").Value = '  -0.59%  '

$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("E13").Value = '  -0.68%  '

$ws.Range("D14").Value = "'3.054.59"
$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("D15").Value = "'58.164.13"
$ws.Range("E15").Value = '  +0.33%  '

$ws.Range("D16").Value = "'20.58"
$ws.Range("E16").Value = '  -1.05%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = "'2.620.36"
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").Value = "'339.99"
$ws.Range("E19").Value = '  +1.92%  '

$ws.Range("D20").Value = "'4.38"
$ws.Range("E20").Value = '  -0.23%  '

$ws.Range("D21").Value = "'10.34"
$ws.Range("E21").Value = '  +0.19%  '

$ws.Range("D22").Value = "'6.41"
$ws.Range("E22").Value = '  +2.74%  '

$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").Value = "'65.99"
$ws.Range("E24").Value = '  +3.33%  '

$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("E26").Value = '  -2.54%  '

$ws.Range("D27").Value = "'2.715.76"
$ws.Range("E27").Value = '  -0.50%  '

$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("D29").Value = "'7.05"
$ws.Range("E29").Value = '  -0.18%  '

$ws.Range("D30").Value = "'0.0₃0753"
$ws.Range("E30").Value = '  -3.79%  '

$ws.Range("E31").Value = '  -0.09%  '

$ws.Range("D32").Value = "'6.27"
$ws.Range("E32").Value = '  -4.71%  '

$ws.Range("E33").Value = '  +1.11%  '

$ws.Range("E34").Value = '  +1.26%  '

$ws.Range("D35").Value = "'149.81"
$ws.Range("E35").Value = '  -0.43%  '

$ws.Range("D36").Value = "'4.06"
$ws.Range("E36").Value = '  -0.59%  '

$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("D38").Value = "'0.868"

$ws.Range("D39").Value = "'0.850"
$ws.Range("E39").Value = '  +1.34%  '

$ws.Range("E40").Value = '  +2.77%  '

$ws.Range("D41").Value = "'36.12"
$ws.Range("E41").Value = '  -0.97%  '

$ws.Range("D42").Value = "'3.56"
$ws.Range("E42").Value = '  -0.76%  '

$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = '  -0.36%  '

$ws.Range("D44").Value = "'274.50"
$ws.Range("E44").Value = '  +2.64%  '

$ws.Range("E45").Value = '  +0.27%  '

$ws.Range("E46").Value = '  -0.54%  '

$ws.Range("E47").Value = '  +0.40%  '

$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("D49").Value = "'18.83"
$ws.Range("E49").Value = '  -1.18%  '

$ws.Range("D50").Value = "'19.18"
$ws.Range("E50").Value = '  +5.74%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = "'4.72"
$ws.Range("E51").Value = '  +2.52%  '
